# Задание 3.4.2, 3.4.3, 3.5.1
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Статистика по годам"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Статистика по годам")

$years = @(
    @(2003, 1366,  0,      1983, 0),
    @(2004, 1488,  0,      7833, 0),
    @(2005, 13331, 0,      16022, 0),
    @(2006, 1522,  0,      33321, 0),
    @(2007, 5604,  27500,  53562, 2),
    @(2008, 27478, 22000,  75070, 2),
    @(2009, 37548, 0,      52889, 0),
    @(2010, 40958, 35000,  93494, 5),
    @(2011, 42359, 33666,  142458, 12),
    @(2012, 44540, 28800,  173897, 18),
    @(2013, 46218, 38050,  234019, 21),
    @(2014, 48482, 29681,  259571, 19),
    @(2015, 50654, 28653,  284763, 27),
    @(2016, 58261, 34538,  332460, 28),
    @(2017, 61724, 35071,  391464, 67),
    @(2018, 65563, 46296,  517670, 71),
    @(2019, 78212, 51657,  535956, 73),
    @(2020, 90537, 40704,  489472, 79),
    @(2021, 105356, 61594, 287915, 48),
    @(2022, 124935, 47928, 91142, 10)
)

# Extend formatting (style) for the newly-added rows 18:21 by copying the
# format of the last existing data row (17) down, before writing values.
$ws1.Range("A17:E17").Copy() | Out-Null
$ws1.Range("A18:E21").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

for ($i = 0; $i -lt $years.Length; $i++) {
    $r = $i + 2
    $row = $years[$i]
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
}

# ---------------------------------------------------------------------------
# Sheet 2: "Статистика по городам"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Статистика по городам")

$cities = @(
    @("Алматы",          179001, "Москва",           0.4917),
    @("Москва",          70218,  "Санкт-Петербург",  0.1493),
    @("Санкт-Петербург",  61230, "Минск",            0.0598),
    @("Новосибирск",      57592, "Киев",             0.0474),
    @("Екатеринбург",     54972, "Новосибирск",      0.0348),
    @("Казань",           50049, "Нижний Новгород",  0.0316),
    @("Краснодар",        47800, "Алматы",           0.029),
    @("Челябинск",        46337, "Воронеж",          0.0274),
    @("Нижний Новгород",  44775, "Казань",           0.0273),
    @("Пермь",            44542, "Воронеж",          0.0141)
)

for ($i = 0; $i -lt $cities.Length; $i++) {
    $r = $i + 2
    $row = $cities[$i]
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 4).Value = $row[2]
    $ws2.Cells.Item($r, 5).Value = $row[3]
}
